# Update the "timestamp" column (H) for data rows 2-51 with the new
# refreshed timestamp value, mirroring an automated data-ingestion run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2025-03-01 03:43:57"

for ($row = 2; $row -le 51; $row++) {
    $ws.Cells.Item($row, 8).Value = $newTimestamp
}
